$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2026-02-16 14:14:26"
$ws.Range("D3").Value = "2026-02-16 14:14:27"
$ws.Range("D4").Value = "2026-02-16 14:14:27"
$ws.Range("D5").Value = "2026-02-16 14:14:27"
$ws.Range("D6").Value = "2026-02-16 14:14:27"
$ws.Range("D7").Value = "2026-02-16 14:14:27"
